$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E (rows 2-51) to Text format so values such as
# "1.000" or "0.000008676" are preserved exactly as strings instead of
# being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

# Apply the cell value changes described by the diff
$ws.Range("D2").Value = "27.676.10"
$ws.Range("E2").Value = "  -0.65%  "
$ws.Range("D3").Value = "1.898.02"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "310.65"
$ws.Range("E5").Value = "  -0.79%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "0.5271"
$ws.Range("E7").Value = "  +6.05%  "
$ws.Range("D8").Value = "0.3798"
$ws.Range("E8").Value = "  -0.61%  "
$ws.Range("D9").Value = "0.07247"
$ws.Range("E9").Value = "  -1.24%  "
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "0.9028"
$ws.Range("E10").Value = "  -0.81%  "
$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").Value = "21.09"
$ws.Range("E11").Value = "  +0.80%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.912.74"
$ws.Range("E12").Value = "  -1.80%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "0.07634"
$ws.Range("E13").Value = "  +0.00%  "
$ws.Range("D14").Value = "5.439"
$ws.Range("E14").Value = "  -0.91%  "
$ws.Range("D15").Value = "91.68"
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").Value = "0.000008676"
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("D18").Value = "14.33"
$ws.Range("E18").Value = "  -1.38%  "
$ws.Range("D19").Value = "0.9996"
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").Value = "27.712.40"
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("D21").Value = "5.145"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").Value = "2.132.68"
$ws.Range("E22").Value = "  -2.81%  "
$ws.Range("D24").Value = "6.593"
$ws.Range("E24").Value = "  -0.68%  "
$ws.Range("D25").Value = "153.45"
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("D26").Value = "1.861"
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("D27").Value = "18.26"
$ws.Range("E27").Value = "  -0.85%  "
$ws.Range("D28").Value = "2.188"
$ws.Range("E28").Value = "  -2.01%  "
$ws.Range("D29").Value = "114.23"
$ws.Range("E29").Value = "  -0.88%  "
$ws.Range("D30").Value = "4.830"
$ws.Range("E30").Value = "  -2.45%  "
$ws.Range("D31").Value = "4.822"
$ws.Range("E31").Value = "  +3.83%  "
$ws.Range("D32").Value = "0.09147"
$ws.Range("E32").Value = "  +1.96%  "
$ws.Range("D33").Value = "0.05267"
$ws.Range("E33").Value = "  -0.55%  "
$ws.Range("D34").Value = "3.144"
$ws.Range("E34").Value = "  -1.57%  "
$ws.Range("E35").Value = "  -1.46%  "
$ws.Range("D36").Value = "0.7725"
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").Value = "0.02084"
$ws.Range("E37").Value = "  +0.94%  "
$ws.Range("D38").Value = "2.569"
$ws.Range("E38").Value = "  -0.28%  "
$ws.Range("D39").Value = "3.077"
$ws.Range("E39").Value = "  +2.22%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "1.090"
$ws.Range("E40").Value = "  -0.91%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.5551"
$ws.Range("E41").Value = "  +0.46%  "
$ws.Range("D42").Value = "6.703"
$ws.Range("E42").Value = "  -4.16%  "
$ws.Range("D43").Value = "117.40"
$ws.Range("E43").Value = "  +5.25%  "
$ws.Range("D44").Value = "8.728"
$ws.Range("E44").Value = "  +1.93%  "
$ws.Range("D45").Value = "0.1511"
$ws.Range("E45").Value = "  -0.87%  "
$ws.Range("D46").Value = "0.4804"
$ws.Range("E46").Value = "  +0.14%  "
$ws.Range("D47").Value = "10.43"
$ws.Range("E47").Value = "  -2.07%  "
$ws.Range("D48").Value = "0.9997"
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("D49").Value = "1.594"
$ws.Range("E49").Value = "  -2.83%  "
$ws.Range("D50").Value = "66.26"
$ws.Range("E50").Value = "  -2.01%  "
$ws.Range("D51").Value = "37.05"
$ws.Range("E51").Value = "  -0.04%  "
